$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.242.83'
$ws.Range('E2').Value = '  -3.06%  '
$ws.Range('D3').Value = '1.728.74'
$ws.Range('E3').Value = '  -3.79%  '
$ws.Range('D4').Value = '''1.006'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''322.27'
$ws.Range('E5').Value = '  -4.14%  '
$ws.Range('D6').Value = '''1.003'
$ws.Range('D7').Value = '''0.4219'
$ws.Range('E7').Value = '  -8.24%  '
$ws.Range('D8').Value = '''0.3568'
$ws.Range('E8').Value = '  -3.39%  '
$ws.Range('D9').Value = '''44.89'
$ws.Range('E9').Value = '  -0.74%  '
$ws.Range('D10').Value = '''0.07385'
$ws.Range('E10').Value = '  -3.31%  '
$ws.Range('E11').Value = '  -3.56%  '
$ws.Range('D12').Value = '''1.004'
$ws.Range('E12').Value = '  +0.08%  '
$ws.Range('E13').Value = '  -5.04%  '
$ws.Range('D14').Value = '''6.040'
$ws.Range('E14').Value = '  -4.60%  '
$ws.Range('D15').Value = '''7.095'
$ws.Range('E15').Value = '  -3.90%  '
$ws.Range('D16').Value = '1.731.27'
$ws.Range('E16').Value = '  -3.70%  '
$ws.Range('D17').Value = '''0.00001053'
$ws.Range('E17').Value = '  -3.38%  '
$ws.Range('D18').Value = '''86.30'
$ws.Range('E18').Value = '  +4.77%  '
$ws.Range('D19').Value = '''0.05932'
$ws.Range('E19').Value = '  -11.66%  '
$ws.Range('E20').Value = '  +0.11%  '
$ws.Range('D21').Value = '''16.71'
$ws.Range('E21').Value = '  -3.93%  '
$ws.Range('D22').Value = '''6.046'
$ws.Range('E22').Value = '  -5.52%  '
$ws.Range('D23').Value = '''0.5250'
$ws.Range('E23').Value = '  -4.73%  '
$ws.Range('D24').Value = '27.293.88'
$ws.Range('E24').Value = '  -2.90%  '
$ws.Range('D25').Value = '''11.28'
$ws.Range('E25').Value = '  -4.90%  '
$ws.Range('E26').Value = '  -1.33%  '
$ws.Range('D27').Value = '''20.02'
$ws.Range('E27').Value = '  -3.49%  '
$ws.Range('D28').Value = '''2.336'
$ws.Range('E28').Value = '  -1.94%  '
$ws.Range('D29').Value = '''147.94'
$ws.Range('E29').Value = '  -2.25%  '
$ws.Range('D30').Value = '1.925.78'
$ws.Range('B31').Value = 'BitcoinCash'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D31').Value = '''125.50'
$ws.Range('E31').Value = '  -5.99%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').Value = '''1.196'
$ws.Range('E32').Value = '  -4.69%  '
$ws.Range('D33').Value = '''0.09067'
$ws.Range('E33').Value = '  -6.06%  '
$ws.Range('D34').Value = '''5.572'
$ws.Range('E34').Value = '  -5.65%  '
$ws.Range('D35').Value = '''3.556'
$ws.Range('E35').Value = '  -12.29%  '
$ws.Range('D36').Value = '''12.59'
$ws.Range('E36').Value = '  +3.59%  '
$ws.Range('D37').Value = '''0.2129'
$ws.Range('E37').Value = '  -4.00%  '
$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').Value = '''5.023'
$ws.Range('E38').Value = '  -4.28%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = '''0.06041'
$ws.Range('E39').Value = '  -4.49%  '
$ws.Range('D40').Value = '''0.02235'
$ws.Range('E40').Value = '  -5.87%  '
$ws.Range('D41').Value = '''0.6326'
$ws.Range('E41').Value = '  -5.58%  '
$ws.Range('D42').Value = '''1.180'
$ws.Range('E42').Value = '  -4.47%  '
$ws.Range('D43').Value = '''1.002'
$ws.Range('E43').Value = '  +0.09%  '
$ws.Range('D44').Value = '''1.411'
$ws.Range('E44').Value = '  -6.24%  '
$ws.Range('D45').Value = '''7.867'
$ws.Range('E45').Value = '  -2.47%  '
$ws.Range('D46').Value = '''13.47'
$ws.Range('E46').Value = '  -4.66%  '
$ws.Range('D48').Value = '''0.5788'
$ws.Range('E48').Value = '  -5.80%  '
$ws.Range('D49').Value = '''123.97'
$ws.Range('E49').Value = '  -4.68%  '
$ws.Range('D50').Value = '''1.935'
$ws.Range('E50').Value = '  -5.75%  '
$ws.Range('D51').Value = '''0.06809'
$ws.Range('E51').Value = '  -4.53%  '
